$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text format for the rows whose price text looks like a plain
# number (e.g. "0.617"), so Excel keeps them as literal text instead of converting them
# to numeric values (this mirrors how the existing data -- e.g. "34.887.05" -- is stored).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.970.58"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.820.41"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "230.43"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "40.19"
$ws.Range("E8").Value = "  -6.04%  "
$ws.Range("D9").Value = "0.323"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "0.0990"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "2.081.77"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "11.30"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").Value = "0.670"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "1.816.43"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "4.61"
$ws.Range("E16").Value = "  -1.82%  "
$ws.Range("D17").Value = "34.931.96"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "69.73"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "0.0₃0785"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("D20").Value = "240.92"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "12.05"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "2.27"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "173.43"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").Value = "7.81"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "0.124"
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("D28").Value = "17.35"
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "3.96"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  +12.48%  "
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").Value = "0.692"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "93.06"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "1.36"
$ws.Range("E38").Value = "  +7.27%  "
$ws.Range("D39").Value = "1.339.24"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "0.979"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").Value = "14.51"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  +2.46%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "1.998.47"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "0.0664"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").Value = "97.41"
$ws.Range("E51").Value = "  -3.33%  "
